$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$data = @{
    2  = @(1430.8, 1393, 1410.85, 1412.45, 121, 1397.15)
    3  = @(1232.45, 1204.35, 1217, 1218.65, 41, 1210.5)
    4  = @(48950, 48423.8, 48505, 48501.35, 27, 48943.25)
    5  = @(118.3, 115.9, 116.35, 116.5, 945, 116.45)
    6  = @(828.5, 818.55, 819.05, 820.35, 136, 820.45)
    7  = @(714.5, 684.7, 705.1, 705.85, 315, 689.45)
    8  = @(1123.2, 1102.4, 1107.75, 1106.45, 347, 1122.6)
    9  = @(1054.8, 1035, 1041.5, 1043.45, 42, 1052.55)
    10 = @(22864.3, 22721.75, 22731, 22743.8, 70, 22845)
    11 = @(2921.9, 2883, 2886.45, 2888, 135, 2912.65)
    12 = @(834.85, 821.55, 827.05, 824.75, 326, 831.6)
    13 = @(1090.75, 1070.95, 1071.9, 1074.5, 36, 1088.35)
    14 = @(954, 942.3, 946.3, 947.7, 282, 944.25)
    15 = @(176.5, 173.9, 174.5, 174.65, 831, 175.85)
    16 = @(3842, 3807.3, 3808.05, 3810.3, 36, 3838.95)
    17 = @(3425, 3378.6, 3388.5, 3387.15, 18, 3410)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 2).Value = $vals[0]
    $ws.Cells.Item($row, 3).Value = $vals[1]
    $ws.Cells.Item($row, 4).Value = $vals[2]
    $ws.Cells.Item($row, 5).Value = $vals[3]
    $ws.Cells.Item($row, 6).Value = $vals[4]
    $ws.Cells.Item($row, 7).Value = $vals[5]
}
